$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell A6 held "[bill_payment_on]"; rename it to "[bill_payment_type]".
$ws.Range("A6").Value = "[bill_payment_type]"

# Move the active selection to the edited cell, matching the authored change.
$ws.Range("A6").Select()
